$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the dated block (rows 321-322),
# pushing the existing rows 321-332 down to 323-334.
$ws.Rows("321:322").Insert()

# New row 321: Camote, 2a (guarda), week of 2021-09-09
$ws.Range("A321").Value = 10
$ws.Range("B321").Value = "Vega Modelo de Temuco"
$ws.Range("C321").Value = "La Araucanía"
$ws.Range("D321").Value = 44448
$ws.Range("E321").Value = 9
$ws.Range("F321").Value = 100112045
$ws.Range("G321").Value = "Zapallo"
$ws.Range("H321").Value = "Camote"
$ws.Range("I321").Value = "2a (guarda)"
$ws.Range("J321").Value = 300
$ws.Range("K321").Value = 800
$ws.Range("L321").Value = 800
$ws.Range("M321").Value = 800
$ws.Range("N321").Value = "$/kilo (volumen en unidades)"
$ws.Range("O321").Value = "Región de O'Higgins"
$ws.Range("P321").Value = 800
$ws.Range("Q321").Value = 1
$ws.Range("R321").Value = "Hortaliza"

# New row 322: Paine, 1a (guarda), week of 2021-09-09
$ws.Range("A322").Value = 10
$ws.Range("B322").Value = "Vega Modelo de Temuco"
$ws.Range("C322").Value = "La Araucanía"
$ws.Range("D322").Value = 44448
$ws.Range("E322").Value = 9
$ws.Range("F322").Value = 100112045
$ws.Range("G322").Value = "Zapallo"
$ws.Range("H322").Value = "Paine"
$ws.Range("I322").Value = "1a (guarda)"
$ws.Range("J322").Value = 2200
$ws.Range("K322").Value = 400
$ws.Range("L322").Value = 500
$ws.Range("M322").Value = 445
$ws.Range("N322").Value = "$/kilo (volumen en unidades)"
$ws.Range("O322").Value = "Región de O'Higgins"
$ws.Range("P322").Value = 445
$ws.Range("Q322").Value = 1
$ws.Range("R322").Value = "Hortaliza"
